$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Pijn_waar (bodyParts)" -> "Pijn_waar" (row 21, column A)
$ws.Range("A21").Value = "Pijn_waar"

# 2. Insert a new row at position 70 (shifts existing rows 70-75 down to 71-76)
#    with the Medication_content translation entry.
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = "Medicatie_inhoud (open)"
$ws.Range("B70").Value = "Medication_content"
$ws.Range("C70").Value = "?"

# 3. Append a new row 77 with the Start_questionnaire_week translation entry.
$ws.Range("A77").Value = "start_vragenlijsten_week"
$ws.Range("B77").Value = "Start_questionnaire_week"
$ws.Range("C77").Value = "?"

# Keep selection/view similar to the committed workbook (top-left cell A64, selection B77).
$ws.Range("B77").Select()
